$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Brasil -> Brasil
$ws.Range("B12").Value = 108620
$ws.Range("C12").Value = 354
$ws.Range("E12").Value = 55438
$ws.Range("G12").Value = 24
$ws.Range("H12").Value = 7367

# Row 24: Mexico -> Mexico
$ws.Range("B24").Value = 24905
$ws.Range("C24").Value = 1434
$ws.Range("E24").Value = 9187
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = 2271

# Row 40: Corea del Sur -> Corea del Sur
$ws.Range("B40").Value = 10804
$ws.Range("C40").Value = 3
$ws.Range("D40").Value = 9283
$ws.Range("E40").Value = 1267
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = 254

# Row 80: Cuba -> Bolivia
$ws.Range("A80").Value = "Bolivia"
$ws.Range("B80").Value = 1681
$ws.Range("C80").Value = 87
$ws.Range("D80").Value = 174
$ws.Range("E80").Value = 1425
$ws.Range("F80").Value = 3
$ws.Range("G80").Value = 6
$ws.Range("H80").Value = 82

# Row 81: Bulgaria -> Cuba
$ws.Range("A81").Value = "Cuba"
$ws.Range("B81").Value = 1668
$ws.Range("D81").Value = 876
$ws.Range("E81").Value = 723
$ws.Range("F81").Value = 9
$ws.Range("H81").Value = 69

# Row 82: Bolivia -> Bulgaria
$ws.Range("A82").Value = "Bulgaria"
$ws.Range("B82").Value = 1652
$ws.Range("D82").Value = 321
$ws.Range("E82").Value = 1253
$ws.Range("F82").Value = 37
$ws.Range("H82").Value = 78

# Row 84: Nueva Zelanda -> Nueva Zelanda
$ws.Range("B84").Value = 1486
$ws.Range("D84").Value = 1302
$ws.Range("E84").Value = 164

# Row 90: Republica de Yibuti -> Honduras
$ws.Range("A90").Value = "Honduras"
$ws.Range("B90").Value = 1178
$ws.Range("C90").Value = 123
$ws.Range("D90").Value = 122
$ws.Range("E90").Value = 973
$ws.Range("F90").Value = 10
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 83

# Row 91: Honduras -> Republica de Yibuti
$ws.Range("A91").Value = "Republica de Yibuti"
$ws.Range("B91").Value = 1116
$ws.Range("D91").Value = 713
$ws.Range("E91").Value = 401
$ws.Range("F91").Value = 0
$ws.Range("H91").Value = 2

# Row 104: Crucero -> Guatemala
$ws.Range("A104").Value = "Guatemala"
$ws.Range("B104").Value = 730
$ws.Range("C104").Value = 27
$ws.Range("D104").Value = 79
$ws.Range("E104").Value = 632
$ws.Range("F104").Value = 5
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 19

# Row 105: Guatemala -> Crucero
$ws.Range("A105").Value = "Crucero"
$ws.Range("B105").Value = 712
$ws.Range("D105").Value = 645
$ws.Range("E105").Value = 54
$ws.Range("F105").Value = 4
$ws.Range("H105").Value = 13

# Row 123: Guinea-Bisau -> Paraguay
$ws.Range("A123").Value = "Paraguay"
$ws.Range("B123").Value = 415
$ws.Range("C123").Value = 19
$ws.Range("D123").Value = 130
$ws.Range("E123").Value = 275
$ws.Range("F123").Value = 7
$ws.Range("H123").Value = 10

# Row 124: Paraguay -> Guinea-Bisau
$ws.Range("A124").Value = "Guinea-Bisau"
$ws.Range("B124").Value = 413
$ws.Range("D124").Value = 19
$ws.Range("E124").Value = 393
$ws.Range("F124").Value = 0
$ws.Range("H124").Value = 1

# Row 190: Belice -> Belice
$ws.Range("D190").Value = 14
$ws.Range("E190").Value = 2

# Row 198: San Cristobal y Nieves -> Burundi
$ws.Range("A198").Value = "Burundi"
$ws.Range("D198").Value = 7
$ws.Range("H198").Value = 1

# Row 199: Burundi -> San Cristobal y Nieves
$ws.Range("A199").Value = "San Cristobal y Nieves"
$ws.Range("D199").Value = 8
$ws.Range("H199").Value = 0

# Update timestamp
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 04:03"
